# template_kelas.xlsx — add a "PETUNJUK" (instructions) block next to the
# Kelas input so users know the expected class-name format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D to fit the new instructional text.
$ws.Columns.Item(4).ColumnWidth = 86.91

# New "PETUNJUK" (instructions) heading next to the Kelas row, bold to match
# the other header cells.
$ws.Range("D2").Value = "PETUNJUK"
$ws.Range("D2").Font.Bold = $true

# Explanation + full example underneath it.
$ws.Range("D3").Value = "Format yang benar: XI (Menggunakan romawi) Nama kelas + No kelas (Contoh: OTOMOTIF 2)"
$ws.Range("D4").Value = "Contoh lengkap: XI OTOMOTIF 2"

# Update the example hint in A2 to use the corrected "OTOMOTIF" casing.
$ws.Range("A2").Value = "(Contoh: XI OTOMOTIF 2)"
